# Updating the Forecast Portfolio
# The ENTSO-E forecast window shifted forward by 16 days (the timestamps in
# column A all move from 2025-12-14/15 to 2025-12-30/31), and the shape of
# the forecasted portfolio flows (columns B-N) shifted along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift every timestamp in column A (rows 2-101) forward by 16 days.
for ($r = 2; $r -le 101; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value2 = $cur + 16
}

# 2) The old "15.22" bump (rows 6-9, cols B/L/N) is gone in the refreshed
#    forecast - zero it out.
$ws.Range("B6:B9").Value2 = 0
$ws.Range("L6:L9").Value2 = 0
$ws.Range("N6:N9").Value2 = 0

# 3) The old "29.08" bump (rows 10-13, cols C/M) and its mirrored deficit in
#    N are gone too.
$ws.Range("C10:C13").Value2 = 0
$ws.Range("M10:M13").Value2 = 0
$ws.Range("N10:N13").Value2 = 0

# 4) Same for the old "391.48" bump (rows 14-17, cols C/M/N).
$ws.Range("C14:C17").Value2 = 0
$ws.Range("M14:M17").Value2 = 0
$ws.Range("N14:N17").Value2 = 0

# 5) The old "267.5" bump (rows 30-33, cols D/L/N) is also gone.
$ws.Range("D30:D33").Value2 = 0
$ws.Range("L30:L33").Value2 = 0
$ws.Range("N30:N33").Value2 = 0

# 6) A new "675.82" bump now appears at rows 38-41, cols B/L/N.
$ws.Range("B38:B41").Value2 = 675.8199999999999
$ws.Range("L38:L41").Value2 = 675.8199999999999
$ws.Range("N38:N41").Value2 = 675.8199999999999
